# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45206 (2023-10-07) to 45208 (2023-10-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 398
$firstRow = 2

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45206) {
        $cell.Value = 45208
    }
}
